$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns I and J - copy formatting from H1 (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..35: values for columns I and J
$data = @(
    @(6, 7),
    @(7, 8),
    @(7, 8),
    @(5, 6),
    @(1, 4),
    @(1, 7),
    @(1, 5),
    @(1, 4),
    @(4, 5),
    @(6, 7),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(5, 6),
    @(6, 7),
    @(6, 8),
    @(1, 7),
    @(1, 7),
    @(1, 4),
    @(6, 8),
    @(5, 9),
    @(3, 3),
    @(5, 6),
    @(7, 9),
    @(1, 3),
    @(3, 4),
    @(1, 2)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
